# The task "notify if file was changed by another program / - or another
# instance of this program / - reload or overwrite" (Id 14) was finished:
# remove it from the "Active" sheet and re-add it as a completed ("Done")
# item at the top of the "Inactive" sheet, stamped with the date it was
# finished (3/5/2018).

$wb = $excel.ActiveWorkbook

$active = $wb.Worksheets.Item("Active")
$inactive = $wb.Worksheets.Item("Inactive")

# Locate the row to move by its Id (column A) so this isn't dependent on
# the row always being exactly row 4.
$taskId = 14
$title = "notify if file was changed by another program`n- or another instance of this program`n- reload or overwrite"
$category = "Bug"
$created = "12/15/2017"
$doneDate = "3/5/2018"

$rowIndex = -1
$lastRow = $active.Cells.Item($active.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    if ($active.Cells.Item($r, 1).Value2 -eq $taskId) {
        $rowIndex = $r
        break
    }
}

if ($rowIndex -ne -1) {
    # Remove the finished task from the Active sheet, shifting rows up.
    $active.Rows.Item($rowIndex).Delete()
}

# Insert a fresh row right under the header of the Inactive sheet and give
# it a plain (unbolded) style so it matches the rest of the data rows.
$inactive.Rows.Item(2).Insert()
$newRow = $inactive.Range("A2:F2")
$newRow.Style = "Normal"

$inactive.Range("A2").Value = $taskId
$inactive.Range("B2").Value = $title
$inactive.Range("C2").Value = "Done"
$inactive.Range("D2").Value = $category
$inactive.Range("E2").Value = "'" + $created
$inactive.Range("F2").Value = "'" + $doneDate
